$d = $word.ActiveDocument

$pairs = @(
  @("212×2=424", "424×8=3392"),
  @("470×4=1880", "431×5=2155"),
  @("695×9=6255", "527×6=3162"),
  @("358×8=2864", "810×6=4860"),
  @("134×4=536", "163×6=978"),
  @("318×3=954", "270×6=1620"),
  @("576×2=1152", "729×4=2916"),
  @("867×8=6936", "627×4=2508"),
  @("376×6=2256", "945×4=3780"),
  @("858×6=5148", "818×2=1636"),
  @("906×5=4530", "762×2=1524"),
  @("709×4=2836", "506×8=4048"),
  @("752×3=2256", "233×9=2097"),
  @("974×3=2922", "942×2=1884"),
  @("179×9=1611", "549×2=1098"),
  @("157×6=942", "469×5=2345"),
  @("272×8=2176", "517×7=3619"),
  @("735×9=6615", "847×5=4235"),
  @("832×4=3328", "447×7=3129"),
  @("964×6=5784", "219×6=1314"),
  @("707×4=2828", "782×4=3128"),
  @("879×6=5274", "751×5=3755"),
  @("380×7=2660", "329×6=1974"),
  @("712×8=5696", "895×2=1790"),
  @("415×8=3320", "987×4=3948")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $rng = $d.Content
    $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
